$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from AC1 to the new header cells AD1:AF1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-53)
$ws.Range("AD2:AD53").Value = 64
$ws.Range("AE2:AE53").Value = 98
$ws.Range("AF2:AF53").Value = 0
